$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: title + link update
$ws.Range("D6").Value = "[Object Detection(객체 검출)] YOLO v1 : You Only Look Once"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Object-Detection%EA%B0%9D%EC%B2%B4-%EA%B2%80%EC%B6%9C-%EB%94%A5%EB%9F%AC%EB%8B%9D-%EC%95%8C%EA%B3%A0%EB%A6%AC%EC%A6%98-history-%EB%B0%8F-%EC%9B%90%EB%A6%AC"

# Row 26: title only
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 37: title + link (uid 1970 -> 1968)
$ws.Range("D37").Value = "[Paper Review] LEARNING TO REMEMBER PATTERNS: PATTERN  MATCHING MEMORY NETWORKS FOR TRAFFIC FORECASTING"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1968&mod=document&pageid=1"

# Row 51: title + link
$ws.Range("D51").Value = "[vscode] 전체 프로젝트에서 어떤 단어를 검색하려면?"
$ws.Range("E51").Value = "https://bskyvision.com/1019"
